$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.794.89'
$ws.Range("E2").Value = '  +4.01%  '
$ws.Range("D3").Value = '2.772.14'
$ws.Range("E3").Value = '  +4.55%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '343.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '115.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.548'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.576'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.56%  '
$ws.Range("E11").Value = '  +4.54%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.131'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.17%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.95'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("D15").Value = '3.212.65'
$ws.Range("E15").Value = '  +4.72%  '
$ws.Range("D16").Value = '2.766.33'
$ws.Range("E16").Value = '  +4.25%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.879'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.13%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '51.633.55'
$ws.Range("E18").Value = '  +3.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.51%  '
$ws.Range("D22").Value = '0.0₃0979'
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '275.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.97%  '
$ws.Range("E25").Value = '  +7.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.17%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("E29").Value = '  +1.01%  '
$ws.Range("E30").Value = '  +1.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.57'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.04'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("E33").Value = '  +3.06%  '
$ws.Range("E34").Value = '  -0.58%  '
$ws.Range("E35").Value = '  +0.16%  '
$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.96'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.37%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.09'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.95'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.22'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0380'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.37%  '
$ws.Range("E41").Value = '  +24.14%  '
$ws.Range("E42").Value = '  +2.66%  '
$ws.Range("E43").Value = '  +3.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '126.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '23.26'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.93%  '
$ws.Range("D46").Value = '2.063.90'
$ws.Range("E46").Value = '  -0.26%  '
$ws.Range("E47").Value = '  -0.82%  '
$ws.Range("E49").Value = '  +4.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.889'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +13.25%  '
$ws.Range("E51").Value = '  -1.01%  '
